# Adds two new worksheets, "CaseDetailStat" and "CaseDetailStat_Message",
# after the existing "StatOutput_Message" sheet. These mirror the existing
# StatOutput / StatOutput_Message sheets but cover the case-level-detail
# query used for the ICDC/CTDC case-detail switch.

$wb = $excel.ActiveWorkbook

# --- common message-block text (already present verbatim in the workbook,
#     reused here so the new sheets share the existing string table entries)
$neo4jUrlLabel = "Neo4j_URL:"
$neo4jUrlValue = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userNameLabel = "User_name:"
$userNameValue = "neo4j"
$pwdLabel = "PWD:"
$pwdValue = "icdcDBneo4j0"
$cypherLabel = "Cypher:"
$caseDetailCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN ['COTC007B'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"
$caseDetailStatCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN ['COTC007B']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$outputLabel = "Output:"
$outputValue = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC01_Canine_Filter_Study-COTB_Neo4jData.xlsx"

# --- 1) CaseDetailStat : a small 4-column stat summary sheet, mirrors StatOutput
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$caseDetailStat = $wb.Worksheets.Add($null, $lastSheet)
$caseDetailStat.Name = "CaseDetailStat"

$caseDetailStat.Range("A1").Value = "number_of_files"
$caseDetailStat.Range("B1").Value = "number_of_sample"
$caseDetailStat.Range("C1").Value = "number_of_cases"
$caseDetailStat.Range("D1").Value = "number_of_study"

$statRow = $caseDetailStat.Range("A2:D2")
$statRow.NumberFormat = "@"
$caseDetailStat.Range("A2").Value = "0"
$caseDetailStat.Range("B2").Value = "0"
$caseDetailStat.Range("C2").Value = "84"
$caseDetailStat.Range("D2").Value = "1"
$statRow.ClearFormats()

# --- 2) CaseDetailStat_Message : connection/query message log, mirrors StatOutput_Message
#     but with a third repeated block for the new case-detail-stat query.
$caseDetailMessage = $wb.Worksheets.Add($null, $caseDetailStat)
$caseDetailMessage.Name = "CaseDetailStat_Message"

$block1 = @($neo4jUrlLabel, $neo4jUrlValue, $userNameLabel, $userNameValue, $pwdLabel, $pwdValue, $cypherLabel, $caseDetailCypher, $outputLabel, $outputValue)
$block2 = @($neo4jUrlLabel, $neo4jUrlValue, $userNameLabel, $userNameValue, $pwdLabel, $pwdValue, $cypherLabel, $caseDetailStatCypher, $outputLabel, $outputValue)

$row = 1
foreach ($val in $block1) {
    $caseDetailMessage.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
foreach ($val in $block2) {
    $caseDetailMessage.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
foreach ($val in $block2) {
    $caseDetailMessage.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
